# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.272.63'
$ws.Range('E2').Value = '  -4.65%  '
$ws.Range('D3').Value = '2.238.20'
$ws.Range('E3').Value = '  -5.44%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.587'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.70%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.562'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.22'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('E12').Value = '  -9.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.06%  '
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '2.572.05'
$ws.Range('E15').Value = '  -5.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.864'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -11.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.74%  '
$ws.Range('D18').Value = '2.233.74'
$ws.Range('E18').Value = '  -5.74%  '
$ws.Range('D19').Value = '43.151.97'
$ws.Range('E19').Value = '  -4.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').Value = '0.0₃0962'
$ws.Range('E21').Value = '  -9.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.19%  '
$ws.Range('E24').Value = '  -10.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '238.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.74%  '
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -9.97%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -11.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.79%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.07%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0874'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -10.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.24%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.94'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.122'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.65%  '
$ws.Range('E41').Value = '  -10.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0324'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.54%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '1.798.06'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.206'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '77.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -14.70%  '
